$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.0981586584325309
$ws.Range("H2").Value = 0.08276122181566332
$ws.Range("I2").Value = 455125.3006419801
$ws.Range("J2").Value = 164948.65052199
$ws.Range("L2").Value = 164948.65052199
$ws.Range("M2").Value = 620073.9511639702
$ws.Range("N2").Value = 10164152.6588
$ws.Range("O2").Value = 9756411.7287
$ws.Range("P2").Value = 0.01622847039582581
$ws.Range("Q2").Value = 0.01690669224595841

$ws.Range("D3").Value = 87
$ws.Range("E3").Value = 0.8446601941747572
$ws.Range("F3").Value = 0.8446601941747572
$ws.Range("G3").Value = 0.09729244173422845
$ws.Range("H3").Value = 0.08217905272696967
$ws.Range("I3").Value = 469976.0611579026
$ws.Range("J3").Value = 169775.6107859813
$ws.Range("L3").Value = 169775.6107859813
$ws.Range("M3").Value = 639751.6719438838
$ws.Range("N3").Value = 10441590.754864
$ws.Range("O3").Value = 10034217.596861
$ws.Range("P3").Value = 0.01625955419741909
$ws.Range("Q3").Value = 0.01691966604741481

$ws.Range("G4").Value = 0.09766085713146683
$ws.Range("H4").Value = 0.08169706317728474
$ws.Range("I4").Value = 500624.5054121671
$ws.Range("J4").Value = 178013.4603193245
$ws.Range("L4").Value = 178013.4603193245
$ws.Range("M4").Value = 678637.9657314916
$ws.Range("N4").Value = 10968786.42890992
$ws.Range("O4").Value = 10560442.07616683
$ws.Range("P4").Value = 0.01622909348021789
$ws.Range("Q4").Value = 0.01685662958381935

$ws.Range("C5").Value = 105
$ws.Range("E5").Value = 0.8285714285714286
$ws.Range("G5").Value = 0.09805477375316646
$ws.Range("H5").Value = 0.08124538396690935
$ws.Range("I5").Value = 515808.63208648
$ws.Range("J5").Value = 183031.5598848782
$ws.Range("L5").Value = 183031.5598848782
$ws.Range("M5").Value = 698840.1919713583
$ws.Range("N5").Value = 11245515.25037722
$ws.Range("O5").Value = 10834820.56705184
$ws.Range("P5").Value = 0.01627596031037693
$ws.Range("Q5").Value = 0.0168929018023121

$ws.Range("G6").Value = 0.09525740397532713
$ws.Range("H6").Value = 0.08087892790357966
$ws.Range("I6").Value = 535448.0544955599
$ws.Range("J6").Value = 189790.0884046672
$ws.Range("L6").Value = 189790.0884046672
$ws.Range("M6").Value = 725238.1429002271
$ws.Range("N6").Value = 11615833.42568854
$ws.Range("O6").Value = 11201367.90186339
$ws.Range("P6").Value = 0.01633891271072676
$ws.Range("Q6").Value = 0.01694347423166905

